$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Row 3 becomes a duplicate of row 2 (same username/password shown),
# picking up row 2's cell formatting (bordered, hyperlink style) in the
# process, then the text is corrected to the ineuron credentials.
$ws.Range("A2:B2").Copy()
$ws.Range("A3:B3").PasteSpecial(-4122)
$ws.Range("A3").Value = "ineuron@ineuron.ai"
$ws.Range("B3").Value = "ineuron"

# Selection ends up covering the whole row that was just edited.
$ws.Range("A3:B3").Select()
